{"js": "// Replace the division problems' text in the table cells per the commit diff.\n// Each old expression is unique in the document, so a simple search+replace\n// (matching the whole cell text, case-sensitive, whole word off since \"\u00f7\"/\"=\"\n// are not word characters) is sufficient and unambiguous.\nconst replacements = [\n  [\"691\u00f77=\", \"104\u00f72=\"],\n  [\"825\u00f75=\", \"294\u00f74=\"],\n  [\"964\u00f79=\", \"132\u00f77=\"],\n  [\"211\u00f77=\", \"802\u00f76=\"],\n  [\"971\u00f75=\", \"382\u00f75=\"],\n  [\"515\u00f78=\", \"643\u00f77=\"],\n  [\"671\u00f78=\", \"598\u00f74=\"],\n  [\"918\u00f73=\", \"297\u00f75=\"],\n  [\"165\u00f77=\", \"399\u00f73=\"],\n  [\"920\u00f79=\", \"645\u00f77=\"],\n  [\"522\u00f76=\", \"507\u00f77=\"],\n  [\"140\u00f78=\", \"718\u00f77=\"],\n  [\"409\u00f77=\", \"587\u00f73=\"],\n  [\"686\u00f74=\", \"991\u00f78=\"],\n  [\"518\u00f72=\", \"281\u00f78=\"],\n  [\"390\u00f74=\", \"469\u00f72=\"],\n  [\"939\u00f73=\", \"627\u00f73=\"],\n  [\"515\u00f76=\", \"885\u00f76=\"],\n  [\"311\u00f72=\", \"442\u00f72=\"],\n  [\"382\u00f74=\", \"564\u00f76=\"],\n  [\"866\u00f74=\", \"164\u00f75=\"],\n  [\"966\u00f79=\", \"883\u00f76=\"],\n  [\"585\u00f73=\", \"748\u00f74=\"],\n  [\"886\u00f77=\", \"633\u00f74=\"],\n  [\"298\u00f74=\", \"899\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems' text in the table cells per the commit diff.\n# Each \"old\" expression occurs exactly once in the document, so Find/Replace\n# with MatchCase + whole-text matching is unambiguous for every pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"691\u00f77=\", \"104\u00f72=\"),\n    @(\"825\u00f75=\", \"294\u00f74=\"),\n    @(\"964\u00f79=\", \"132\u00f77=\"),\n    @(\"211\u00f77=\", \"802\u00f76=\"),\n    @(\"971\u00f75=\", \"382\u00f75=\"),\n    @(\"515\u00f78=\", \"643\u00f77=\"),\n    @(\"671\u00f78=\", \"598\u00f74=\"),\n    @(\"918\u00f73=\", \"297\u00f75=\"),\n    @(\"165\u00f77=\", \"399\u00f73=\"),\n    @(\"920\u00f79=\", \"645\u00f77=\"),\n    @(\"522\u00f76=\", \"507\u00f77=\"),\n    @(\"140\u00f78=\", \"718\u00f77=\"),\n    @(\"409\u00f77=\", \"587\u00f73=\"),\n    @(\"686\u00f74=\", \"991\u00f78=\"),\n    @(\"518\u00f72=\", \"281\u00f78=\"),\n    @(\"390\u00f74=\", \"469\u00f72=\"),\n    @(\"939\u00f73=\", \"627\u00f73=\"),\n    @(\"515\u00f76=\", \"885\u00f76=\"),\n    @(\"311\u00f72=\", \"442\u00f72=\"),\n    @(\"382\u00f74=\", \"564\u00f76=\"),\n    @(\"866\u00f74=\", \"164\u00f75=\"),\n    @(\"966\u00f79=\", \"883\u00f76=\"),\n    @(\"585\u00f73=\", \"748\u00f74=\"),\n    @(\"886\u00f77=\", \"633\u00f74=\"),\n    @(\"298\u00f74=\", \"899\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(=1 wdFindContinue), Format, ReplaceWith,\n    # Replace(=2 wdReplaceAll)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
